# Apply the change described by the diff: split merged paragraphs into
# separate sentences joined by manual line breaks (<w:br/>) instead of
# being glued together with no space between sentences.

$d = $word.ActiveDocument

function Replace-WithBreak($findText, $replaceText) {
    $find = $d.Content.Find
    $find.ClearFormatting()
    $ok = $find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)
    if (-not $ok) {
        throw "Find/Replace failed for: $findText"
    }
}

# --- Paragraph 1 (Objetivos, Portuguese) ---
Replace-WithBreak `
    "Modelo Artigo ou Modelo Produto.No Modelo Artigo" `
    "Modelo Artigo ou Modelo Produto.^lNo Modelo Artigo"

Replace-WithBreak `
    "resultados esperados.No Modelo Produto" `
    "resultados esperados.^lNo Modelo Produto"

# --- Paragraph 2 (Objetivos, English/italic) ---
Replace-WithBreak `
    "Article Model or Product Model.In the Article Model" `
    "Article Model or Product Model.^lIn the Article Model"

Replace-WithBreak `
    "expected results.In the Product Model" `
    "expected results.^lIn the Product Model"

# --- Avaliação / Método ---
Replace-WithBreak `
    "indicados pelo aluno.Modelo Produto" `
    "indicados pelo aluno.^lModelo Produto"

Replace-WithBreak `
    "Engenharia Ambiental.Em ambos modelos" `
    "Engenharia Ambiental.^lEm ambos modelos"

# --- Avaliação / Critério ---
Replace-WithBreak `
    "de igual pesoFica sob responsabilidade" `
    "de igual peso^lFica sob responsabilidade"
